$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The match result for row 2 (Turkey vs Italy) changed: previously Turkey (A) won 3-0,
# now Italy (B) won 0-3. Update score1/score2 and the Awon/Bwon indicator columns;
# the LossA/LossB/Loss formula columns will recalculate automatically.
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 3
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 1

$excel.Calculate()

# Update the selected cell to match the saved view state.
$ws.Range("R23").Select()
